$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 9 new rows right after the header row (before the current row 2),
# shifting all existing data rows down by 9.
$ws.Range("A2:A10").EntireRow.Insert()
# The inserted rows pick up formatting copied from the row above (the
# header); clear it so the new data rows stay unstyled like the rest of
# the data rows.
$ws.Range("A2:C10").ClearFormats()

# Fill the newly inserted rows (2-10) with the new accelerometer readings.
$newRows = @(
    @(0.2728629112243653, -0.005532175302505469, 0.17908151820302),
    @(-0.1289855480194095, -0.0322514295578003, 0.2833916515111924),
    @(-0.2541066646575928, -0.01140816211700431, 0.3126487381756305),
    @(-0.1588943481445311, -0.05039391517639156, 0.1764491081237791),
    @(0.003789019584655856, -0.07255983948707578, 0.1795819453895092),
    @(0.0466578006744384, -0.03053182363510123, 0.2458087503910065),
    @(-0.2323491334915166, -0.003547763824462882, 0.2952604919672013),
    @(-0.07431058883666952, 0.007624650001525833, 0.2976654559373855),
    @(0.08380470275878904, 0.005895948410034098, 0.3545163981616498)
)

$r = 2
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r++
}

# Append a new row (31) at the end with another accelerometer reading.
$ws.Cells.Item(31, 1).Value = -0.5065834045410164
$ws.Cells.Item(31, 2).Value = 1.347906202077864
$ws.Cells.Item(31, 3).Value = 1.883451831340787
